$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("A9").Value = "Notes:"
$ws.Range("A9").Font.Bold = $true

$ws.Range("A10").Value = "This variable affects what portion of newly sold vehicles qualify for Vehicle Battery Production subsidy,"
$ws.Range("A11").Value = "if relevant. We assume that the U.S. has enough battery manufacturing capacity to supply"
$ws.Range("A12").Value = "100% of vehicle battery demand domestically. All U.S. state models should use the U.S. vaues, since the "
$ws.Range("A13").Value = "battery manufacturing production tax credit applies to all U.S. manufactured batteries."

$ws.Range("A14").Select()
